$wb = $excel.ActiveWorkbook

# Update "展览" sheet (F3: 110 -> 111, F4: 61 -> 62)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 111
$ws1.Range("F4").Value = 62

# Update "全部类型" sheet (F3: 110 -> 111, F4: 61 -> 62)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 111
$ws4.Range("F4").Value = 62
